$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.051.06"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.572.37"
$ws.Range("E3").Value = "  -1.62%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'572.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "'142.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "2.575.62"
$ws.Range("E9").Value = "  -2.12%  "
$ws.Range("D10").Value = "'6.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("E12").Value = "  +11.48%  "
$ws.Range("D13").Value = "'0.343"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").Value = "3.025.90"
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").Value = "59.093.27"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").Value = "'22.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "2.583.86"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "'4.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").Value = "'335.62"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").Value = "'10.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "'6.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").Value = "'64.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.43%  "
$ws.Range("E25").Value = "  +8.01%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("E27").Value = "  -1.60%  "
$ws.Range("D28").Value = "'7.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "0.0₃0778"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "'1.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'159.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("D33").Value = "'6.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "'18.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "'4.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.43%  "
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("D38").Value = "'0.874"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.68%  "
$ws.Range("D39").Value = "'37.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "'294.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("E42").Value = "  +1.12%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'131.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.77%  "
$ws.Range("E45").Value = "  +1.92%  "
$ws.Range("D46").Value = "'0.593"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "'0.0535"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "'10.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").Value = "'19.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.17%  "
$ws.Range("E50").Value = "  +1.64%  "
$ws.Range("D51").Value = "1.947.02"
$ws.Range("E51").Value = "  -0.34%  "
